# Ajuste leitura de tabelas
# Realinha os nomes da coluna A e os totais das colunas B:J conforme a leitura
# corrigida das tabelas de origem.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Novos nomes (coluna A), linhas 2 a 24, na ordem correta apos o ajuste
$names = @(
    "Tamires Teixeira",
    "Natalia Farias",
    "Daniela Fernandes",
    "Vania Fagundes",
    "Luana Umpierre",
    "Nicolas Silva",
    "Adriana Hunhoff",
    "Carlla Bo",
    "Dominique Daudt",
    "Josue Alos",
    "Igor Martins",
    "Alexia Pereira",
    "Gabriel Winter",
    "Aline Castro",
    "Julio Acauan",
    "Daniel Machado",
    "Brenda Pereira",
    "Amanda Bernardes",
    "Eduarda Santos",
    "Michele Mattidorff",
    "Jonathan Cardoso",
    "Gabriel Wolff",
    "Brenda Fossa"
)

# Novos valores numericos (colunas B..J), linhas 2 a 25 (25 = linha de totais)
$data = @(
    @(4,11,73,4,19,12,1,1,110),
    @(7,13,66,10,18,0,0,0,94),
    @(3,1,3,31,41,0,0,1,76),
    @(3,0,0,43,17,0,0,4,64),
    @(0,0,0,12,40,1,0,0,53),
    @(1,0,0,15,24,0,0,2,41),
    @(0,0,0,10,17,0,0,0,27),
    @(0,0,0,13,11,0,0,0,24),
    @(1,0,0,11,11,0,0,1,23),
    @(0,0,0,4,5,0,0,0,9),
    @(0,0,0,4,0,0,0,0,4),
    @(0,0,0,1,1,0,0,0,2),
    @(0,1,0,2,0,0,0,0,2),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(19,26,142,160,204,13,1,9,529)
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
}

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $col = $c + 2
        $ws.Cells.Item($row, $col).Value = $vals[$c]
    }
}
